# chore: update Sheets via scheduled runner
# Refresh market-derived profit figures (currentAveragePrice / NQ / HQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ, columns H:N) for
# the affected Leve rows across the Anima_Profits sheets. Values are plain
# numeric literals pulled from an external market snapshot, not formulas,
# so each changed cell is written directly; cells that drop out of the
# source snapshot are cleared so the <c> element is omitted on save, and
# cells newly present in the snapshot are written so a new <c> element
# appears.
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 4
$ws.Range("H4").Value = 216.66667
$ws.Range("I4").Value = 216.66667
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 216.66667
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -102.66667
$ws.Range("N4").Value = $null
# Row 22
$ws.Range("H22").Value = 1500
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = $null
# Row 30
$ws.Range("H30").Value = 2850.6667
$ws.Range("I30").Value = 836
$ws.Range("K30").Value = 2508
$ws.Range("M30").Value = -2407
# Row 129
$ws.Range("H129").Value = 1486.5122
$ws.Range("I129").Value = 537.7692
$ws.Range("J129").Value = 1927
$ws.Range("K129").Value = 1613.3076
$ws.Range("L129").Value = 5781
$ws.Range("M129").Value = 3386.6924
$ws.Range("N129").Value = -15781
# Row 132
$ws.Range("H132").Value = 3015.25
$ws.Range("I132").Value = 2882.4783
$ws.Range("J132").Value = 4033.1667
$ws.Range("K132").Value = 8647.4349
$ws.Range("L132").Value = 12099.5001
$ws.Range("M132").Value = -6117.4349
$ws.Range("N132").Value = -17159.5001
# Row 137
$ws.Range("H137").Value = 1310.8788
$ws.Range("I137").Value = 1294.091
$ws.Range("J137").Value = 1327.6666
$ws.Range("K137").Value = 3882.273
$ws.Range("L137").Value = 3982.9998
$ws.Range("M137").Value = -1332.273
$ws.Range("N137").Value = -9082.9998
# Row 138
$ws.Range("H138").Value = 1213.36
$ws.Range("I138").Value = 537.04254
$ws.Range("J138").Value = 1813.1132
$ws.Range("K138").Value = 1611.12762
$ws.Range("L138").Value = 5439.3396
$ws.Range("M138").Value = 3528.87238
$ws.Range("N138").Value = -15719.3396
# Row 141
$ws.Range("H141").Value = 2531.3584
$ws.Range("I141").Value = 867.3261
$ws.Range("J141").Value = 13466.429
$ws.Range("K141").Value = 2601.9783
$ws.Range("L141").Value = 40399.287
$ws.Range("M141").Value = 2578.0217
$ws.Range("N141").Value = -50759.287

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 24
$ws.Range("H24").Value = 45335
$ws.Range("J24").Value = 45335
$ws.Range("L24").Value = 45335
$ws.Range("N24").Value = -46083
# Row 32
$ws.Range("H32").Value = 1813.06
$ws.Range("I32").Value = 1672.1333
$ws.Range("J32").Value = 3081.4
$ws.Range("K32").Value = 1672.1333
$ws.Range("L32").Value = 3081.4
$ws.Range("M32").Value = -1385.1333
$ws.Range("N32").Value = -3655.4
# Row 61
$ws.Range("H61").Value = 6411818
$ws.Range("I61").Value = 6945802.5
$ws.Range("J61").Value = 3999.75
$ws.Range("K61").Value = 6945802.5
$ws.Range("L61").Value = 3999.75
$ws.Range("M61").Value = -6945590.5
$ws.Range("N61").Value = -4423.75
# Row 74
$ws.Range("H74").Value = 1592.5424
$ws.Range("I74").Value = 887.80646
$ws.Range("J74").Value = 2372.7856
$ws.Range("K74").Value = 887.80646
$ws.Range("L74").Value = 2372.7856
$ws.Range("M74").Value = -13.80646000000002
$ws.Range("N74").Value = -4120.7856
# Row 77
$ws.Range("H77").Value = 1592.5424
$ws.Range("I77").Value = 887.80646
$ws.Range("J77").Value = 2372.7856
$ws.Range("K77").Value = 4439.0323
$ws.Range("L77").Value = 11863.928
$ws.Range("M77").Value = -71.03229999999985
$ws.Range("N77").Value = -20599.928
# Row 100
$ws.Range("H100").Value = 45335
$ws.Range("J100").Value = 45335
$ws.Range("L100").Value = 45335
$ws.Range("N100").Value = -47499
# Row 132
$ws.Range("H132").Value = 2950.25
$ws.Range("I132").Value = 2597.4285
$ws.Range("J132").Value = 3900.1538
$ws.Range("K132").Value = 7792.2855
$ws.Range("L132").Value = 11700.4614
$ws.Range("M132").Value = -5262.2855
$ws.Range("N132").Value = -16760.4614
# Row 136
$ws.Range("H136").Value = 6411818
$ws.Range("I136").Value = 6945802.5
$ws.Range("J136").Value = 3999.75
$ws.Range("K136").Value = 20837407.5
$ws.Range("L136").Value = 11999.25
$ws.Range("M136").Value = -20834857.5
$ws.Range("N136").Value = -17099.25

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 7
$ws.Range("H7").Value = 32778230
$ws.Range("I7").Value = 65555556
$ws.Range("J7").Value = 904
$ws.Range("K7").Value = 65555556
$ws.Range("L7").Value = 904
$ws.Range("M7").Value = -65555443
$ws.Range("N7").Value = -1130
# Row 134
$ws.Range("H134").Value = 2082.5715
$ws.Range("I134").Value = 2115.9375
$ws.Range("J134").Value = 1726.6666
$ws.Range("K134").Value = 6347.8125
$ws.Range("L134").Value = 5179.9998
$ws.Range("M134").Value = -3812.8125
$ws.Range("N134").Value = -10249.9998

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 200
$ws.Range("I7").Value = 0
$ws.Range("K7").Value = 0
$ws.Range("M7").Value = $null
# Row 31
$ws.Range("H31").Value = 3830.0396
$ws.Range("I31").Value = 1266.6177
$ws.Range("J31").Value = 5905.1904
$ws.Range("K31").Value = 1266.6177
$ws.Range("L31").Value = 5905.1904
$ws.Range("M31").Value = -971.6177
$ws.Range("N31").Value = -6495.1904
# Row 34
$ws.Range("H34").Value = 3830.0396
$ws.Range("I34").Value = 1266.6177
$ws.Range("J34").Value = 5905.1904
$ws.Range("K34").Value = 1266.6177
$ws.Range("L34").Value = 5905.1904
$ws.Range("M34").Value = -1064.6177
$ws.Range("N34").Value = -6309.1904
# Row 58
$ws.Range("H58").Value = 1248.909
$ws.Range("I58").Value = 1048.0769
$ws.Range("J58").Value = 1539
$ws.Range("K58").Value = 1048.0769
$ws.Range("L58").Value = 1539
$ws.Range("M58").Value = -845.0769
$ws.Range("N58").Value = -1945
# Row 86
$ws.Range("H86").Value = 2668
$ws.Range("I86").Value = 3000
$ws.Range("K86").Value = 3000
$ws.Range("M86").Value = -1877
# Row 89
$ws.Range("H89").Value = 2668
$ws.Range("I89").Value = 3000
$ws.Range("K89").Value = 12500
$ws.Range("M89").Value = -9384
# Row 132
$ws.Range("H132").Value = 3206517.2
$ws.Range("I132").Value = 1137.85
$ws.Range("J132").Value = 13891115
$ws.Range("K132").Value = 3413.55
$ws.Range("L132").Value = 41673345
$ws.Range("M132").Value = -883.5499999999997
$ws.Range("N132").Value = -41678405
# Row 134
$ws.Range("H134").Value = 3498.2
$ws.Range("I134").Value = 3331.186
$ws.Range("J134").Value = 4524.143
$ws.Range("K134").Value = 9993.558000000001
$ws.Range("L134").Value = 13572.429
$ws.Range("M134").Value = -7458.558000000001
$ws.Range("N134").Value = -18642.429
# Row 136
$ws.Range("H136").Value = 1248.909
$ws.Range("I136").Value = 1048.0769
$ws.Range("J136").Value = 1539
$ws.Range("K136").Value = 3144.2307
$ws.Range("L136").Value = 4617
$ws.Range("M136").Value = -594.2307000000001
$ws.Range("N136").Value = -9717

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 1351.9678
$ws.Range("I5").Value = 445.05554
$ws.Range("J5").Value = 2607.6924
$ws.Range("K5").Value = 1335.16662
$ws.Range("L5").Value = 7823.0772
$ws.Range("M5").Value = -1223.16662
$ws.Range("N5").Value = -8047.0772
# Row 107
$ws.Range("H107").Value = 20833646
$ws.Range("J107").Value = 58823876
$ws.Range("L107").Value = 176471628
$ws.Range("N107").Value = -176475468
# Row 114
$ws.Range("H114").Value = 1354.75
$ws.Range("J114").Value = 1706.3334
$ws.Range("L114").Value = 5119.0002
$ws.Range("N114").Value = -11627.0002
# Row 131
$ws.Range("H131").Value = 2959.6724
$ws.Range("I131").Value = 418.42856
$ws.Range("J131").Value = 3768.25
$ws.Range("K131").Value = 1255.28568
$ws.Range("L131").Value = 11304.75
$ws.Range("M131").Value = 3784.71432
$ws.Range("N131").Value = -21384.75
# Row 135
$ws.Range("H135").Value = 1351.9678
$ws.Range("I135").Value = 445.05554
$ws.Range("J135").Value = 2607.6924
$ws.Range("K135").Value = 4005.49986
$ws.Range("L135").Value = 23469.2316
$ws.Range("M135").Value = -1470.49986
$ws.Range("N135").Value = -28539.2316

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 100
$ws.Range("H100").Value = 3796.2222
$ws.Range("I100").Value = 3793.2
$ws.Range("J100").Value = 3800
$ws.Range("K100").Value = 3793.2
$ws.Range("L100").Value = 3800
$ws.Range("M100").Value = -3252.2
$ws.Range("N100").Value = -4882
# Row 115
$ws.Range("H115").Value = 40000
$ws.Range("J115").Value = 40000
$ws.Range("L115").Value = 40000
$ws.Range("N115").Value = -42350
# Row 132
$ws.Range("H132").Value = 2153.4358
$ws.Range("I132").Value = 1966.2667
$ws.Range("J132").Value = 2777.3333
$ws.Range("K132").Value = 5898.800099999999
$ws.Range("L132").Value = 8331.999899999999
$ws.Range("M132").Value = -3368.800099999999
$ws.Range("N132").Value = -13391.9999
# Row 136
$ws.Range("H136").Value = 3334592.2
$ws.Range("I136").Value = 1115.0303
$ws.Range("J136").Value = 9805460
$ws.Range("K136").Value = 3345.0909
$ws.Range("L136").Value = 29416380
$ws.Range("M136").Value = -795.0908999999997
$ws.Range("N136").Value = -29421480

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 8
$ws.Range("H8").Value = 100
$ws.Range("I8").Value = 100
$ws.Range("K8").Value = 100
$ws.Range("M8").Value = 40
# Row 132
$ws.Range("H132").Value = 4420653
$ws.Range("I132").Value = 1574.2106
$ws.Range("J132").Value = 10417975
$ws.Range("K132").Value = 4722.6318
$ws.Range("L132").Value = 31253925
$ws.Range("M132").Value = -2192.6318
$ws.Range("N132").Value = -31258985
# Row 136
$ws.Range("H136").Value = 1313.39
$ws.Range("I136").Value = 1230.3026
$ws.Range("J136").Value = 1576.5
$ws.Range("K136").Value = 3690.9078
$ws.Range("L136").Value = 4729.5
$ws.Range("M136").Value = -1140.9078
